# Auto-generated Excel COM-interop script to update '想去人数' (F column)
# values across all 4 worksheets per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 2496
$ws.Range("F10").Value = 5024
$ws.Range("F11").Value = 6241
$ws.Range("F12").Value = 901
$ws.Range("F14").Value = 1407
$ws.Range("F15").Value = 1349
$ws.Range("F16").Value = 561
$ws.Range("F17").Value = 6832
$ws.Range("F18").Value = 381
$ws.Range("F19").Value = 32
$ws.Range("F21").Value = 63
$ws.Range("F22").Value = 4598
$ws.Range("F23").Value = 377
$ws.Range("F24").Value = 32
$ws.Range("F26").Value = 2240
$ws.Range("F28").Value = 420
$ws.Range("F29").Value = 1134
$ws.Range("F30").Value = 192
$ws.Range("F31").Value = 80
$ws.Range("F33").Value = 136
$ws.Range("F35").Value = 1255
$ws.Range("F36").Value = 1963
$ws.Range("F37").Value = 202
$ws.Range("F38").Value = 485
$ws.Range("F40").Value = 1334
$ws.Range("F41").Value = 585
$ws.Range("F43").Value = 6
$ws.Range("F44").Value = 1072
$ws.Range("F45").Value = 1689
$ws.Range("F46").Value = 36
$ws.Range("F49").Value = 62

$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 353
$ws.Range("F13").Value = 246
$ws.Range("F16").Value = 164
$ws.Range("F21").Value = 17
$ws.Range("F27").Value = 256
$ws.Range("F40").Value = 3

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 488
$ws.Range("F6").Value = 1631
$ws.Range("F7").Value = 529
$ws.Range("F8").Value = 1215
$ws.Range("F10").Value = 1705
$ws.Range("F11").Value = 2068
$ws.Range("F12").Value = 532
$ws.Range("F13").Value = 436

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 488
$ws.Range("F3").Value = 1631
$ws.Range("F6").Value = 2496
$ws.Range("F8").Value = 1215
$ws.Range("F10").Value = 5024
$ws.Range("F11").Value = 532
$ws.Range("F14").Value = 901
$ws.Range("F17").Value = 1407
$ws.Range("F18").Value = 1349
$ws.Range("F19").Value = 561
$ws.Range("F20").Value = 6832
$ws.Range("F21").Value = 381
$ws.Range("F22").Value = 436
$ws.Range("F25").Value = 32
$ws.Range("F26").Value = 4598
$ws.Range("F27").Value = 2240
$ws.Range("F28").Value = 420
$ws.Range("F29").Value = 1134
$ws.Range("F30").Value = 192
$ws.Range("F31").Value = 80
$ws.Range("F33").Value = 246
$ws.Range("F34").Value = 136
$ws.Range("F36").Value = 1255
$ws.Range("F37").Value = 1963
$ws.Range("F38").Value = 202
$ws.Range("F39").Value = 485
$ws.Range("F41").Value = 17
$ws.Range("F42").Value = 1334
$ws.Range("F44").Value = 6
$ws.Range("F46").Value = 1072
$ws.Range("F47").Value = 1693
$ws.Range("F48").Value = 62
